$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cache a blank, default-formatted cell to use as a "format donor" so that
# forcing text on numeric-looking values does not leave a stray NumberFormat
# style applied to the edited cells (matches original unstyled inlineStr cells).
$fmtDonor = $ws.Range("A1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '300.99'
$fmtDonor.Copy()
$ws.Range('D2').PasteSpecial(-4122)
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-0.01%'
$fmtDonor.Copy()
$ws.Range('E2').PasteSpecial(-4122)
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '3.94%'
$fmtDonor.Copy()
$ws.Range('E3').PasteSpecial(-4122)
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.942'
$fmtDonor.Copy()
$ws.Range('D4').PasteSpecial(-4122)
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-2.91%'
$fmtDonor.Copy()
$ws.Range('E4').PasteSpecial(-4122)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07714'
$fmtDonor.Copy()
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-1.69%'
$fmtDonor.Copy()
$ws.Range('E5').PasteSpecial(-4122)
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.981'
$fmtDonor.Copy()
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-15.41%'
$fmtDonor.Copy()
$ws.Range('E6').PasteSpecial(-4122)
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.832'
$fmtDonor.Copy()
$ws.Range('D7').PasteSpecial(-4122)
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '0.23%'
$fmtDonor.Copy()
$ws.Range('E7').PasteSpecial(-4122)
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9199'
$fmtDonor.Copy()
$ws.Range('D8').PasteSpecial(-4122)
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.22%'
$fmtDonor.Copy()
$ws.Range('E8').PasteSpecial(-4122)
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1754'
$fmtDonor.Copy()
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.09%'
$fmtDonor.Copy()
$ws.Range('E9').PasteSpecial(-4122)
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07753'
$fmtDonor.Copy()
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.77%'
$fmtDonor.Copy()
$ws.Range('E10').PasteSpecial(-4122)
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08614'
$fmtDonor.Copy()
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-6.65%'
$fmtDonor.Copy()
$ws.Range('E11').PasteSpecial(-4122)
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03175'
$fmtDonor.Copy()
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '6.06%'
$fmtDonor.Copy()
$ws.Range('E12').PasteSpecial(-4122)
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1003'
$fmtDonor.Copy()
$ws.Range('D13').PasteSpecial(-4122)
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.04%'
$fmtDonor.Copy()
$ws.Range('E13').PasteSpecial(-4122)
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001510'
$fmtDonor.Copy()
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.50%'
$fmtDonor.Copy()
$ws.Range('E14').PasteSpecial(-4122)
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005914'
$fmtDonor.Copy()
$ws.Range('D15').PasteSpecial(-4122)
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1.99%'
$fmtDonor.Copy()
$ws.Range('E15').PasteSpecial(-4122)
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.458'
$fmtDonor.Copy()
$ws.Range('D16').PasteSpecial(-4122)
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.36%'
$fmtDonor.Copy()
$ws.Range('E16').PasteSpecial(-4122)
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.800'
$fmtDonor.Copy()
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.93%'
$fmtDonor.Copy()
$ws.Range('E17').PasteSpecial(-4122)
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.152'
$fmtDonor.Copy()
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-4.33%'
$fmtDonor.Copy()
$ws.Range('E18').PasteSpecial(-4122)
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3349'
$fmtDonor.Copy()
$ws.Range('D19').PasteSpecial(-4122)
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.37%'
$fmtDonor.Copy()
$ws.Range('E19').PasteSpecial(-4122)
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1326'
$fmtDonor.Copy()
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.13%'
$fmtDonor.Copy()
$ws.Range('E20').PasteSpecial(-4122)
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.280'
$fmtDonor.Copy()
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '5.62%'
$fmtDonor.Copy()
$ws.Range('E21').PasteSpecial(-4122)
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.1991'
$fmtDonor.Copy()
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '11.25%'
$fmtDonor.Copy()
$ws.Range('E22').PasteSpecial(-4122)
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04509'
$fmtDonor.Copy()
$ws.Range('D23').PasteSpecial(-4122)
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.43%'
$fmtDonor.Copy()
$ws.Range('E23').PasteSpecial(-4122)
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001224'
$fmtDonor.Copy()
$ws.Range('D24').PasteSpecial(-4122)
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.99%'
$fmtDonor.Copy()
$ws.Range('E24').PasteSpecial(-4122)
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004419'
$fmtDonor.Copy()
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.97%'
$fmtDonor.Copy()
$ws.Range('E25').PasteSpecial(-4122)
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.11%'
$fmtDonor.Copy()
$ws.Range('E26').PasteSpecial(-4122)
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01706'
$fmtDonor.Copy()
$ws.Range('D39').PasteSpecial(-4122)
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-3.17%'
$fmtDonor.Copy()
$ws.Range('E39').PasteSpecial(-4122)
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04698'
$fmtDonor.Copy()
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-1.17%'
$fmtDonor.Copy()
$ws.Range('E40').PasteSpecial(-4122)
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007487'
$fmtDonor.Copy()
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '5.69%'
$fmtDonor.Copy()
$ws.Range('E41').PasteSpecial(-4122)
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1352'
$fmtDonor.Copy()
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.64%'
$fmtDonor.Copy()
$ws.Range('E42').PasteSpecial(-4122)
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002332'
$fmtDonor.Copy()
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '6.51%'
$fmtDonor.Copy()
$ws.Range('E43').PasteSpecial(-4122)
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01045'
$fmtDonor.Copy()
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '1.12%'
$fmtDonor.Copy()
$ws.Range('E44').PasteSpecial(-4122)
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00006250'
$fmtDonor.Copy()
$ws.Range('D45').PasteSpecial(-4122)
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-0.15%'
$fmtDonor.Copy()
$ws.Range('E45').PasteSpecial(-4122)
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.10%'
$fmtDonor.Copy()
$ws.Range('E46').PasteSpecial(-4122)
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.8204'
$fmtDonor.Copy()
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-28.79%'
$fmtDonor.Copy()
$ws.Range('E47').PasteSpecial(-4122)
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.10%'
$fmtDonor.Copy()
$ws.Range('E49').PasteSpecial(-4122)
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.10%'
$fmtDonor.Copy()
$ws.Range('E50').PasteSpecial(-4122)

$excel.CutCopyMode = $false
